$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6874376666666667
$ws.Range("H2").Value = 2.062313
$ws.Range("I2").Value = 0.2330845252991127
$ws.Range("J2").Value = 0.2330845252991127
$ws.Range("M2").Value = 92.253011
$ws.Range("N2").Value = 276.759033
$ws.Range("O2").Value = 0.2854710184133813
$ws.Range("P2").Value = 0.2854710184133813
$ws.Range("Q2").Value = 63.41819462481434
$ws.Range("R2").Value = 570.763751623329
$ws.Range("S2").Value = 0.06653887681353726
$ws.Range("T2").Value = 0.06653887681353723
$ws.Range("G3").Value = 0.6874376666666667
$ws.Range("H3").Value = 2.062313
$ws.Range("I3").Value = 0.2330845252991127
$ws.Range("J3").Value = 0.2330845252991127
$ws.Range("O3").Value = 0.1080355352256351
$ws.Range("P3").Value = 0.1080355352256351
$ws.Range("Q3").Value = 24.00039989143134
$ws.Range("R3").Value = 216.003599022882
$ws.Range("S3").Value = 0.02518141144350272
$ws.Range("T3").Value = 0.02518141144350272
$ws.Range("G4").Value = 0.6874376666666667
$ws.Range("H4").Value = 2.062313
$ws.Range("I4").Value = 0.2330845252991127
$ws.Range("J4").Value = 0.2330845252991127
$ws.Range("M4").Value = 42.21774566666667
$ws.Range("N4").Value = 126.653237
$ws.Range("O4").Value = 0.1306401029076487
$ws.Range("P4").Value = 0.1306401029076487
$ws.Range("Q4").Value = 29.02206857302012
$ws.Range("R4").Value = 261.198617157181
$ws.Range("S4").Value = 0.03045018637125654
$ws.Range("T4").Value = 0.03045018637125654
$ws.Range("G5").Value = 0.6874376666666667
$ws.Range("H5").Value = 2.062313
$ws.Range("I5").Value = 0.2330845252991127
$ws.Range("J5").Value = 0.2330845252991127
$ws.Range("M5").Value = 13.65158233333333
$ws.Range("N5").Value = 40.954747
$ws.Range("O5").Value = 0.0422439448794879
$ws.Range("P5").Value = 0.0422439448794879
$ws.Range("Q5").Value = 9.384611905534555
$ws.Range("R5").Value = 84.461507149811
$ws.Range("S5").Value = 0.009846409838997322
$ws.Range("T5").Value = 0.009846409838997318
$ws.Range("G6").Value = 0.6874376666666667
$ws.Range("H6").Value = 2.062313
$ws.Range("I6").Value = 0.2330845252991127
$ws.Range("J6").Value = 0.2330845252991127
$ws.Range("M6").Value = 18.17840666666667
$ws.Range("N6").Value = 54.53522
$ws.Range("O6").Value = 0.0562519120841046
$ws.Range("P6").Value = 0.05625191208410459
$ws.Range("Q6").Value = 12.49652146265111
$ws.Range("R6").Value = 112.46869316386
$ws.Range("S6").Value = 0.01311145022529094
$ws.Range("T6").Value = 0.01311145022529094
$ws.Range("G7").Value = 0.6874376666666667
$ws.Range("H7").Value = 2.062313
$ws.Range("I7").Value = 0.2330845252991127
$ws.Range("J7").Value = 0.2330845252991127
$ws.Range("M7").Value = 121.9471053333333
$ws.Range("N7").Value = 365.841316
$ws.Range("O7").Value = 0.3773574864897424
$ws.Range("P7").Value = 0.3773574864897424
$ws.Range("Q7").Value = 83.8310335471009
$ws.Range("R7").Value = 754.479301923908
$ws.Range("S7").Value = 0.08795619060652794
$ws.Range("T7").Value = 0.08795619060652793
$ws.Range("I8").Value = 0.1587189032810992
$ws.Range("J8").Value = 0.1587189032810992
$ws.Range("M8").Value = 92.253011
$ws.Range("N8").Value = 276.759033
$ws.Range("O8").Value = 0.2854710184133813
$ws.Range("P8").Value = 0.2854710184133813
$ws.Range("Q8").Value = 43.18461848121734
$ws.Range("R8").Value = 388.661566330956
$ws.Range("S8").Value = 0.04530964696111037
$ws.Range("T8").Value = 0.04530964696111035
$ws.Range("I9").Value = 0.1587189032810992
$ws.Range("J9").Value = 0.1587189032810992
$ws.Range("O9").Value = 0.1080355352256351
$ws.Range("P9").Value = 0.1080355352256351
$ws.Range("S9").Value = 0.01714728166639936
$ws.Range("T9").Value = 0.01714728166639936
$ws.Range("I10").Value = 0.1587189032810992
$ws.Range("J10").Value = 0.1587189032810992
$ws.Range("M10").Value = 42.21774566666667
$ws.Range("N10").Value = 126.653237
$ws.Range("O10").Value = 0.1306401029076487
$ws.Range("P10").Value = 0.1306401029076487
$ws.Range("Q10").Value = 19.76257706918712
$ws.Range("R10").Value = 177.863193622684
$ws.Range("S10").Value = 0.02073505385803195
$ws.Range("T10").Value = 0.02073505385803195
$ws.Range("I11").Value = 0.1587189032810992
$ws.Range("J11").Value = 0.1587189032810992
$ws.Range("M11").Value = 13.65158233333333
$ws.Range("N11").Value = 40.954747
$ws.Range("O11").Value = 0.0422439448794879
$ws.Range("P11").Value = 0.0422439448794879
$ws.Range("Q11").Value = 6.390451307111555
$ws.Range("R11").Value = 57.514061764004
$ws.Range("S11").Value = 0.006704912601539528
$ws.Range("T11").Value = 0.006704912601539525
$ws.Range("I12").Value = 0.1587189032810992
$ws.Range("J12").Value = 0.1587189032810992
$ws.Range("M12").Value = 18.17840666666667
$ws.Range("N12").Value = 54.53522
$ws.Range("O12").Value = 0.0562519120841046
$ws.Range("P12").Value = 0.05625191208410459
$ws.Range("Q12").Value = 8.509506063671113
$ws.Range("R12").Value = 76.58555457304001
$ws.Range("S12").Value = 0.008928241793453895
$ws.Range("T12").Value = 0.008928241793453891
$ws.Range("I13").Value = 0.1587189032810992
$ws.Range("J13").Value = 0.1587189032810992
$ws.Range("M13").Value = 121.9471053333333
$ws.Range("N13").Value = 365.841316
$ws.Range("O13").Value = 0.3773574864897424
$ws.Range("P13").Value = 0.3773574864897424
$ws.Range("Q13").Value = 57.0847407756569
$ws.Range("R13").Value = 513.7626669809121
$ws.Range("S13").Value = 0.05989376640056412
$ws.Range("T13").Value = 0.05989376640056412
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.758325333333333
$ws.Range("H14").Value = 5.274976
$ws.Range("I14").Value = 0.5961826730104559
$ws.Range("J14").Value = 0.5961826730104558
$ws.Range("M14").Value = 92.253011
$ws.Range("N14").Value = 276.759033
$ws.Range("O14").Value = 0.2854710184133813
$ws.Range("P14").Value = 0.2854710184133813
$ws.Range("Q14").Value = 162.2108063175787
$ws.Range("R14").Value = 1459.897256858208
$ws.Range("S14").Value = 0.1701928748247068
$ws.Range("T14").Value = 0.1701928748247067
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.758325333333333
$ws.Range("H15").Value = 5.274976
$ws.Range("I15").Value = 0.5961826730104559
$ws.Range("J15").Value = 0.5961826730104558
$ws.Range("O15").Value = 0.1080355352256351
$ws.Range("P15").Value = 0.1080355352256351
$ws.Range("Q15").Value = 61.38812751396267
$ws.Range("R15").Value = 552.493147625664
$ws.Range("S15").Value = 0.06440891417093439
$ws.Range("T15").Value = 0.06440891417093438
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.758325333333333
$ws.Range("H16").Value = 5.274976
$ws.Range("I16").Value = 0.5961826730104559
$ws.Range("J16").Value = 0.5961826730104558
$ws.Range("M16").Value = 42.21774566666667
$ws.Range("N16").Value = 126.653237
$ws.Range("O16").Value = 0.1306401029076487
$ws.Range("P16").Value = 0.1306401029076487
$ws.Range("Q16").Value = 74.23253172192356
$ws.Range("R16").Value = 668.0927854973121
$ws.Range("S16").Value = 0.07788536575384306
$ws.Range("T16").Value = 0.07788536575384304
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.758325333333333
$ws.Range("H17").Value = 5.274976
$ws.Range("I17").Value = 0.5961826730104559
$ws.Range("J17").Value = 0.5961826730104558
$ws.Range("M17").Value = 13.65158233333333
$ws.Range("N17").Value = 40.954747
$ws.Range("O17").Value = 0.0422439448794879
$ws.Range("P17").Value = 0.0422439448794879
$ws.Range("Q17").Value = 24.00392305678578
$ws.Range("R17").Value = 216.035307511072
$ws.Range("S17").Value = 0.02518510797675946
$ws.Range("T17").Value = 0.02518510797675945
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1.758325333333333
$ws.Range("H18").Value = 5.274976
$ws.Range("I18").Value = 0.5961826730104559
$ws.Range("J18").Value = 0.5961826730104558
$ws.Range("M18").Value = 18.17840666666667
$ws.Range("N18").Value = 54.53522
$ws.Range("O18").Value = 0.0562519120841046
$ws.Range("P18").Value = 0.05625191208410459
$ws.Range("Q18").Value = 31.96355296163556
$ws.Range("R18").Value = 287.67197665472
$ws.Range("S18").Value = 0.03353641530825065
$ws.Range("T18").Value = 0.03353641530825063
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1.758325333333333
$ws.Range("H19").Value = 5.274976
$ws.Range("I19").Value = 0.5961826730104559
$ws.Range("J19").Value = 0.5961826730104558
$ws.Range("M19").Value = 121.9471053333333
$ws.Range("N19").Value = 365.841316
$ws.Range("O19").Value = 0.3773574864897424
$ws.Range("P19").Value = 0.3773574864897424
$ws.Range("Q19").Value = 214.4226846342684
$ws.Range("R19").Value = 1929.804161708416
$ws.Range("S19").Value = 0.2249739949759616
$ws.Range("T19").Value = 0.2249739949759616
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.01541033333333333
$ws.Range("H20").Value = 0.046231
$ws.Range("I20").Value = 0.005225070437466708
$ws.Range("J20").Value = 0.005225070437466708
$ws.Range("M20").Value = 92.253011
$ws.Range("N20").Value = 276.759033
$ws.Range("O20").Value = 0.2854710184133813
$ws.Range("P20").Value = 0.2854710184133813
$ws.Range("Q20").Value = 1.421649650513667
$ws.Range("R20").Value = 12.794846854623
$ws.Range("S20").Value = 0.001491606179065273
$ws.Range("T20").Value = 0.001491606179065273
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.01541033333333333
$ws.Range("H21").Value = 0.046231
$ws.Range("I21").Value = 0.005225070437466708
$ws.Range("J21").Value = 0.005225070437466708
$ws.Range("O21").Value = 0.1080355352256351
$ws.Range("P21").Value = 0.1080355352256351
$ws.Range("Q21").Value = 0.5380184711926667
$ws.Range("R21").Value = 4.842166240734
$ws.Range("S21").Value = 0.000564493281303359
$ws.Range("T21").Value = 0.000564493281303359
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0.3333333333333333
$ws.Range("G22").Value = 0.01541033333333333
$ws.Range("H22").Value = 0.046231
$ws.Range("I22").Value = 0.005225070437466708
$ws.Range("J22").Value = 0.005225070437466708
$ws.Range("M22").Value = 42.21774566666667
$ws.Range("N22").Value = 126.653237
$ws.Range("O22").Value = 0.1306401029076487
$ws.Range("P22").Value = 0.1306401029076487
$ws.Range("Q22").Value = 0.6505895333052223
$ws.Range("R22").Value = 5.855305799747001
$ws.Range("S22").Value = 0.0006826037396503639
$ws.Range("T22").Value = 0.0006826037396503639
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0.3333333333333333
$ws.Range("G23").Value = 0.01541033333333333
$ws.Range("H23").Value = 0.046231
$ws.Range("I23").Value = 0.005225070437466708
$ws.Range("J23").Value = 0.005225070437466708
$ws.Range("M23").Value = 13.65158233333333
$ws.Range("N23").Value = 40.954747
$ws.Range("O23").Value = 0.0422439448794879
$ws.Range("P23").Value = 0.0422439448794879
$ws.Range("Q23").Value = 0.2103754342841111
$ws.Range("R23").Value = 1.893378908557
$ws.Range("S23").Value = 0.0002207275875517854
$ws.Range("T23").Value = 0.0002207275875517853
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.3333333333333333
$ws.Range("G24").Value = 0.01541033333333333
$ws.Range("H24").Value = 0.046231
$ws.Range("I24").Value = 0.005225070437466708
$ws.Range("J24").Value = 0.005225070437466708
$ws.Range("M24").Value = 18.17840666666667
$ws.Range("N24").Value = 54.53522
$ws.Range("O24").Value = 0.0562519120841046
$ws.Range("P24").Value = 0.05625191208410459
$ws.Range("Q24").Value = 0.2801353062022222
$ws.Range("R24").Value = 2.52121775582
$ws.Range("S24").Value = 0.0002939202028816312
$ws.Range("T24").Value = 0.0002939202028816312
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0.3333333333333333
$ws.Range("G25").Value = 0.01541033333333333
$ws.Range("H25").Value = 0.046231
$ws.Range("I25").Value = 0.005225070437466708
$ws.Range("J25").Value = 0.005225070437466708
$ws.Range("M25").Value = 121.9471053333333
$ws.Range("N25").Value = 365.841316
$ws.Range("O25").Value = 0.3773574864897424
$ws.Range("P25").Value = 0.3773574864897424
$ws.Range("Q25").Value = 1.879245542221778
$ws.Range("R25").Value = 16.913209879996
$ws.Range("S25").Value = 0.001971719447014296
$ws.Range("T25").Value = 0.001971719447014296
$ws.Range("G26").Value = 0.02002233333333333
$ws.Range("H26").Value = 0.060067
$ws.Range("I26").Value = 0.006788827971865474
$ws.Range("J26").Value = 0.006788827971865474
$ws.Range("M26").Value = 92.253011
$ws.Range("N26").Value = 276.759033
$ws.Range("O26").Value = 0.2854710184133813
$ws.Range("P26").Value = 0.2854710184133813
$ws.Range("Q26").Value = 1.847120537245667
$ws.Range("R26").Value = 16.624084835211
$ws.Range("S26").Value = 0.001938013634961687
$ws.Range("T26").Value = 0.001938013634961687
$ws.Range("G27").Value = 0.02002233333333333
$ws.Range("H27").Value = 0.060067
$ws.Range("I27").Value = 0.006788827971865474
$ws.Range("J27").Value = 0.006788827971865474
$ws.Range("O27").Value = 0.1080355352256351
$ws.Range("P27").Value = 0.1080355352256351
$ws.Range("Q27").Value = 0.6990364800486667
$ws.Range("R27").Value = 6.291328320438001
$ws.Range("S27").Value = 0.0007334346634952492
$ws.Range("T27").Value = 0.0007334346634952492
$ws.Range("G28").Value = 0.02002233333333333
$ws.Range("H28").Value = 0.060067
$ws.Range("I28").Value = 0.006788827971865474
$ws.Range("J28").Value = 0.006788827971865474
$ws.Range("M28").Value = 42.21774566666667
$ws.Range("N28").Value = 126.653237
$ws.Range("O28").Value = 0.1306401029076487
$ws.Range("P28").Value = 0.1306401029076487
$ws.Range("Q28").Value = 0.845297776319889
$ws.Range("R28").Value = 7.607679986879002
$ws.Range("S28").Value = 0.0008868931848668299
$ws.Range("T28").Value = 0.0008868931848668299
$ws.Range("G29").Value = 0.02002233333333333
$ws.Range("H29").Value = 0.060067
$ws.Range("I29").Value = 0.006788827971865474
$ws.Range("J29").Value = 0.006788827971865474
$ws.Range("M29").Value = 13.65158233333333
$ws.Range("N29").Value = 40.954747
$ws.Range("O29").Value = 0.0422439448794879
$ws.Range("P29").Value = 0.0422439448794879
$ws.Range("Q29").Value = 0.2733365320054444
$ws.Range("R29").Value = 2.460028788049
$ws.Range("S29").Value = 0.0002867868746398107
$ws.Range("T29").Value = 0.0002867868746398107
$ws.Range("G30").Value = 0.02002233333333333
$ws.Range("H30").Value = 0.060067
$ws.Range("I30").Value = 0.006788827971865474
$ws.Range("J30").Value = 0.006788827971865474
$ws.Range("M30").Value = 18.17840666666667
$ws.Range("N30").Value = 54.53522
$ws.Range("O30").Value = 0.0562519120841046
$ws.Range("P30").Value = 0.05625191208410459
$ws.Range("Q30").Value = 0.3639741177488889
$ws.Range("R30").Value = 3.27576705974
$ws.Range("S30").Value = 0.0003818845542274868
$ws.Range("T30").Value = 0.0003818845542274867
$ws.Range("G31").Value = 0.02002233333333333
$ws.Range("H31").Value = 0.060067
$ws.Range("I31").Value = 0.006788827971865474
$ws.Range("J31").Value = 0.006788827971865474
$ws.Range("M31").Value = 121.9471053333333
$ws.Range("N31").Value = 365.841316
$ws.Range("O31").Value = 0.3773574864897424
$ws.Range("P31").Value = 0.3773574864897424
$ws.Range("Q31").Value = 2.441665592019111
$ws.Range("R31").Value = 21.974990328172
$ws.Range("S31").Value = 0.002561815059674411
$ws.Range("T31").Value = 0.002561815059674411
